$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1002.1053
$ws.Range("I40").Value = 785.4545000000001
$ws.Range("J40").Value = 1300
$ws.Range("K40").Value = 785.4545000000001
$ws.Range("L40").Value = 1300
$ws.Range("M40").Value = -610.4545000000001
$ws.Range("N40").Value = -1650
$ws.Range("H64").Value = 4828.5713
$ws.Range("J64").Value = 5700
$ws.Range("L64").Value = 5700
$ws.Range("N64").Value = -6196
$ws.Range("H67").Value = 4828.5713
$ws.Range("J67").Value = 5700
$ws.Range("L67").Value = 5700
$ws.Range("N67").Value = -7416
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H70").Value = 1261.3529
$ws.Range("I70").Value = 1318.4615
$ws.Range("K70").Value = 3955.3845
$ws.Range("M70").Value = -3685.3845
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H73").Value = 1261.3529
$ws.Range("I73").Value = 1318.4615
$ws.Range("K73").Value = 3955.3845
$ws.Range("M73").Value = -3019.3845
$ws.Range("J112").Value = 1181.6522
$ws.Range("L112").Value = 3544.9566
$ws.Range("N112").Value = -5760.9566
$ws.Range("H137").Value = 91257.16
$ws.Range("I137").Value = 116189.46
$ws.Range("J137").Value = 3994.1
$ws.Range("K137").Value = 348568.38
$ws.Range("L137").Value = 11982.3
$ws.Range("M137").Value = -346018.38
$ws.Range("N137").Value = -17082.3
$ws.Range("H138").Value = 3811.5278
$ws.Range("I138").Value = 5062
$ws.Range("J138").Value = 3655.2188
$ws.Range("K138").Value = 15186
$ws.Range("L138").Value = 10965.6564
$ws.Range("M138").Value = -10046
$ws.Range("N138").Value = -21245.6564

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18088.143
$ws.Range("I32").Value = 12381.25
$ws.Range("K32").Value = 12381.25
$ws.Range("M32").Value = -12094.25
$ws.Range("H61").Value = 21962828
$ws.Range("I61").Value = 35137864
$ws.Range("K61").Value = 35137864
$ws.Range("M61").Value = -35137652
$ws.Range("H74").Value = 23811112
$ws.Range("I74").Value = 35715028
$ws.Range("K74").Value = 35715028
$ws.Range("M74").Value = -35714154
$ws.Range("H77").Value = 23811112
$ws.Range("I77").Value = 35715028
$ws.Range("K77").Value = 178575140
$ws.Range("M77").Value = -178570772
$ws.Range("H110").Value = 1958.2
$ws.Range("I110").Value = 1261.8182
$ws.Range("K110").Value = 1261.8182
$ws.Range("M110").Value = 783.1818000000001
$ws.Range("H132").Value = 14724013
$ws.Range("I132").Value = 21742214
$ws.Range("K132").Value = 65226642
$ws.Range("M132").Value = -65224112
$ws.Range("H135").Value = 50430
$ws.Range("J135").Value = 50430
$ws.Range("L135").Value = 50430
$ws.Range("N135").Value = -60570
$ws.Range("H136").Value = 21962828
$ws.Range("I136").Value = 35137864
$ws.Range("K136").Value = 105413592
$ws.Range("M136").Value = -105411042

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2002.8
$ws.Range("I86").Value = 1859.3125
$ws.Range("J86").Value = 2576.75
$ws.Range("K86").Value = 1859.3125
$ws.Range("L86").Value = 2576.75
$ws.Range("M86").Value = -736.3125
$ws.Range("N86").Value = -4822.75
$ws.Range("H89").Value = 2002.8
$ws.Range("I89").Value = 1859.3125
$ws.Range("J89").Value = 2576.75
$ws.Range("K89").Value = 9296.5625
$ws.Range("L89").Value = 12883.75
$ws.Range("M89").Value = -3680.5625
$ws.Range("N89").Value = -24115.75
$ws.Range("H134").Value = 3124.8794
$ws.Range("I134").Value = 2815.1
$ws.Range("K134").Value = 8445.299999999999
$ws.Range("M134").Value = -5910.299999999999

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6999.93
$ws.Range("I31").Value = 8125
$ws.Range("J31").Value = 6884.5386
$ws.Range("K31").Value = 8125
$ws.Range("L31").Value = 6884.5386
$ws.Range("M31").Value = -7830
$ws.Range("N31").Value = -7474.5386
$ws.Range("H34").Value = 6999.93
$ws.Range("I34").Value = 8125
$ws.Range("J34").Value = 6884.5386
$ws.Range("K34").Value = 8125
$ws.Range("L34").Value = 6884.5386
$ws.Range("M34").Value = -7923
$ws.Range("N34").Value = -7288.5386
$ws.Range("H58").Value = 21959.6
$ws.Range("I58").Value = 2460.5715
$ws.Range("J58").Value = 29542.555
$ws.Range("K58").Value = 2460.5715
$ws.Range("L58").Value = 29542.555
$ws.Range("M58").Value = -2257.5715
$ws.Range("N58").Value = -29948.555
$ws.Range("H122").Value = 1478.381
$ws.Range("I122").Value = 1188.9231
$ws.Range("K122").Value = 3566.7693
$ws.Range("M122").Value = -1116.7693
$ws.Range("H132").Value = 40003880
$ws.Range("I132").Value = 50002900
$ws.Range("J132").Value = 7802.4
$ws.Range("K132").Value = 150008700
$ws.Range("L132").Value = 23407.2
$ws.Range("M132").Value = -150006170
$ws.Range("N132").Value = -28467.2
$ws.Range("H134").Value = 37037964
$ws.Range("I134").Value = 41667564
$ws.Range("J134").Value = 1168
$ws.Range("K134").Value = 125002692
$ws.Range("L134").Value = 3504
$ws.Range("M134").Value = -125000157
$ws.Range("N134").Value = -8574
$ws.Range("H136").Value = 21959.6
$ws.Range("I136").Value = 2460.5715
$ws.Range("J136").Value = 29542.555
$ws.Range("K136").Value = 7381.7145
$ws.Range("L136").Value = 88627.66500000001
$ws.Range("M136").Value = -4831.7145
$ws.Range("N136").Value = -93727.66500000001
$ws.Range("H141").Value = 23140.75
$ws.Range("J141").Value = 23140.75
$ws.Range("L141").Value = 23140.75
$ws.Range("N141").Value = -33500.75

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 8418.833000000001
$ws.Range("J2").Value = 128.57143
$ws.Range("L2").Value = 771.42858
$ws.Range("N2").Value = -997.42858
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H37").Value = 35825716
$ws.Range("J37").Value = 35825716
$ws.Range("L37").Value = 107477148
$ws.Range("N37").Value = -107477372
$ws.Range("H131").Value = 721.5599999999999
$ws.Range("I131").Value = 470
$ws.Range("J131").Value = 724.101
$ws.Range("K131").Value = 1410
$ws.Range("L131").Value = 2172.303
$ws.Range("M131").Value = 3630
$ws.Range("N131").Value = -12252.303
$ws.Range("H132").Value = 751.8
$ws.Range("I132").Value = 751.8
$ws.Range("K132").Value = 6766.2
$ws.Range("M132").Value = -4236.2
$ws.Range("H134").Value = 4500.56
$ws.Range("I134").Value = 3102.7273
$ws.Range("J134").Value = 5598.857
$ws.Range("K134").Value = 9308.1819
$ws.Range("L134").Value = 16796.571
$ws.Range("M134").Value = -4238.1819
$ws.Range("N134").Value = -26936.571

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3977.75
$ws.Range("J97").Value = 3303.6667
$ws.Range("L97").Value = 3303.6667
$ws.Range("N97").Value = -4295.6667
$ws.Range("H105").Value = 30600
$ws.Range("J105").Value = 30600
$ws.Range("L105").Value = 30600
$ws.Range("N105").Value = -37588
$ws.Range("H132").Value = 5317523.5
$ws.Range("I132").Value = 7942947
$ws.Range("K132").Value = 23828841
$ws.Range("M132").Value = -23826311

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2578.2
$ws.Range("I22").Value = 2722.75
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 2722.75
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -2427.75
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 2578.2
$ws.Range("I27").Value = 2722.75
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 2722.75
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -2615.75
$ws.Range("N27").Value = -2214
$ws.Range("H132").Value = 3054
$ws.Range("I132").Value = 1944.6154
$ws.Range("K132").Value = 5833.8462
$ws.Range("M132").Value = -3303.8462
$ws.Range("H136").Value = 3214.2856
$ws.Range("I136").Value = 2916.6667
$ws.Range("K136").Value = 8750.000100000001
$ws.Range("M136").Value = -6200.000100000001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 30535
$ws.Range("J92").Value = 30535
$ws.Range("L92").Value = 30535
$ws.Range("N92").Value = -35527
$ws.Range("H109").Value = 26990
$ws.Range("J109").Value = 26990
$ws.Range("L109").Value = 26990
$ws.Range("N109").Value = -29764
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 23811508
$ws.Range("I132").Value = 38462976
$ws.Range("J132").Value = 2872.25
$ws.Range("K132").Value = 115388928
$ws.Range("L132").Value = 8616.75
$ws.Range("M132").Value = -115386398
$ws.Range("N132").Value = -13676.75
